$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 164
$ws.Range("B4").Value  = 522.2067784721304
$ws.Range("B6").Value  = 3397.251084247763
$ws.Range("B7").Value  = 0.8883134593884953
$ws.Range("B8").Value  = 0.8462854921458476
$ws.Range("B9").Value  = 0.8219161189494576
$ws.Range("B10").Value = 0.8462854921458477
$ws.Range("B11").Value = 0.8706070566086225
$ws.Range("B12").Value = 0.08648754114113126
$ws.Range("B13").Value = 88.43544861111174
$ws.Range("B14").Value = 257.7233886258969
$ws.Range("B15").Value = 1.782275694444131
